# Finland Veikkausliiga workbook update
# The match-odds data rows got re-ordered/re-matched to the correct fixtures.
# For each group below, the full row content in columns B:AC (everything
# except the rank number in column A) is rotated among the listed rows:
#   new(rows[i]) = old(rows[(i+1) mod n])
# i.e. row rows[i] receives the data that used to live in row rows[i+1],
# and the last row in the list wraps around and receives the data that
# used to live in the first row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cycles = @(
    ,@(27, 28)
    ,@(29, 30)
    ,@(44, 45)
    ,@(59, 62, 60, 63)
    ,@(65, 66)
    ,@(71, 73, 72)
    ,@(84, 85)
    ,@(86, 87)
    ,@(102, 103)
    ,@(107, 108)
    ,@(128, 133, 130, 129)
    ,@(137, 138)
    ,@(141, 142)
    ,@(158, 159)
)

foreach ($rows in $cycles) {
    # Snapshot the original B:AC content of every row in this cycle before
    # writing anything, so the rotation uses consistent "before" values.
    $snapshot = @{}
    foreach ($r in $rows) {
        $snapshot[$r] = $ws.Range("B$r`:AC$r").Value2
    }

    $n = $rows.Length
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $rows[$i]
        $srcRow = $rows[($i + 1) % $n]
        $ws.Range("B$destRow`:AC$destRow").Value2 = $snapshot[$srcRow]
    }
}
